# The author fixed a typo in the "돈까스" (pork cutlet) word-problem text:
# the second "bkg" ("돼지고기 bkg으로") should read "akg" ("돼지고기 akg으로"),
# matching the pattern already used by the other similar problems in the
# sheet (e.g. "밀가루 akg으로", "점토 akg으로").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B3")
$text = $target.Value2
$fixed = $text -replace "bkg으로", "akg으로"
$target.Value = $fixed

# Reflect the author's final selection/scroll position in the sheet.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollColumn = 2
